# Apply updates to the teacher's weekly schedule sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 3
$ws.Range("B3").Value = "Circuitos Elétricos 2 - ELT-2A"
$ws.Range("C3").Value = "-"
$ws.Range("D3").Value = "Circuitos Elétricos 2 - MCT-2A"

# Row 4
$ws.Range("C4").Value = "Circuitos Elétricos 2 - ELT-2A"
$ws.Range("D4").Value = "Circuitos Elétricos 2 - MCT-2A"
$ws.Range("E4").Value = "Programação - MCT-2A"

# Row 6
$ws.Range("C6").Value = "Circuitos Elétricos 2 - ELT-2A"
$ws.Range("D6").Value = "-"
$ws.Range("E6").Value = "-"

$wb.Save()
